# Apply stock-count corrections to CryCompanywiseStockReport.
# Two kinds of edits:
#  1. Row pairs that got swapped (B,E,F,G values exchanged between the two rows)
#     - these are duplicate-named items whose stock figures were mixed up.
#  2. Single-row quantity corrections, where the quantity (F) is reduced and
#     the value (G) is recomputed as F * D (rate).
# After the row-level edits, the "Sub Total:" rows and the grand-total rows
# are recomputed as sums, exactly as the report itself computes them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2, $cols) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

# --- 1. Swap pairs (B, E, F, G) ---
$pairs = @(
    @(136,137),
    @(233,234),
    @(246,247),
    @(277,278),
    @(292,293),
    @(294,296),
    @(299,300),
    @(467,468),
    @(472,473),
    @(479,480),
    @(485,486),
    @(591,592),
    @(602,603)
)

foreach ($p in $pairs) {
    Swap-Rows $p[0] $p[1] @("B","E","F","G")
}

# --- 2. Single quantity corrections: set F, recompute G = F * D ---
$singleEdits = @{
    80  = 136
    220 = 16
    328 = 589
    777 = 465
    780 = 134
}

foreach ($row in $singleEdits.Keys) {
    $newQty = $singleEdits[$row]
    $rate = $ws.Range("D$row").Value2
    $ws.Range("F$row").Value2 = $newQty
    $ws.Range("G$row").Value2 = [math]::Round($rate * $newQty, 2)
}

# --- 3. Recompute "Sub Total:" rows that cover the changed data blocks ---
$subtotalBlocks = @{
    114 = @(70,113)
    222 = @(207,221)
    339 = @(276,338)
    781 = @(773,780)
}

foreach ($subRow in $subtotalBlocks.Keys) {
    $range = $subtotalBlocks[$subRow]
    $startRow = $range[0]
    $endRow = $range[1]
    $sum = 0.0
    for ($r = $startRow; $r -le $endRow; $r++) {
        $sum += $ws.Range("G$r").Value2
    }
    $ws.Range("B$subRow").Value2 = [math]::Round($sum, 2)
}

# --- 4. Recompute the grand-total rows (799 = sum of every "Sub Total:" row
#         except itself; 800 mirrors 799). ---
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$grandSum = 0.0
for ($r = 1; $r -le $lastRow; $r++) {
    $label = $ws.Range("A$r").Value2
    if ($label -eq "Sub Total:" -and $r -ne 799) {
        $grandSum += $ws.Range("B$r").Value2
    }
}
$grandSum = [math]::Round($grandSum, 2)

$ws.Range("B799").Value2 = $grandSum
$ws.Range("B800").Value2 = $grandSum
